$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27 (shifts existing rows 27+ down by one)
$ws.Rows.Item(27).Insert()

# Row 32 (previously row 31 before the insert) now uses a shorter model-parameter note
$ws.Range("B32").Value = "{n_estimators=10}"

# New row 27: list of top-10 features used with even fewer features
$ws.Range("A27").Value = 'Title_Mr, "Sex", "Title_Mrs", "Pclass_3", "Title_Miss", "Cabin_NA", "Fare_0", "Age_3", "Age_2", "Embarked_C"'

# Update the active selection like the source edit
[void]$ws.Range("B27").Select()
